# Atualização de bases das ligas, do dia: 07-03-2024 às 23:43
#
# Row 124 (match already played - FTHG/FTAG/FTR and closing odds added/updated)
# Row 125 (new row - match already played)
# Row 126 (new row - upcoming match, only opening odds known)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 124: add result columns (H,I,J) and update closing odds (N..AC)
# ---------------------------------------------------------------------------
$ws.Range("H124").Value = 2
$ws.Range("I124").Value = 3
$ws.Range("J124").Value = "A"

$ws.Range("N124").Value = 3.25
$ws.Range("O124").Value = 3.25
$ws.Range("P124").Value = 2
$ws.Range("R124").Value = 2
$ws.Range("S124").Value = 1.8
$ws.Range("U124").Value = 1.95
$ws.Range("V124").Value = 1.75
$ws.Range("W124").Value = -1
$ws.Range("X124").Value = -1
$ws.Range("Y124").Value = 1
$ws.Range("Z124").Value = -1
$ws.Range("AA124").Value = 0.8
$ws.Range("AB124").Value = 0.95
$ws.Range("AC124").Value = -1

# ---------------------------------------------------------------------------
# Row 125: brand-new match row (already played)
# ---------------------------------------------------------------------------
$ws.Range("A124").Copy()
$ws.Range("A125").PasteSpecial(-4122)
$ws.Range("E124").Copy()
$ws.Range("E125").PasteSpecial(-4122)

$ws.Range("A125").Value = 123
$ws.Range("B125").Value = 7011615
$ws.Range("C125").Value = "Azerbaijan Premier League"
$ws.Range("D125").Value = "Azerbaijan Premier League"
$ws.Range("E125").Value = 45354.5
$ws.Range("F125").Value = "Neftchi Baku"
$ws.Range("G125").Value = "FK Qarabag"
$ws.Range("H125").Value = 1
$ws.Range("I125").Value = 4
$ws.Range("J125").Value = "A"
$ws.Range("K125").Value = 4.75
$ws.Range("L125").Value = 3.8
$ws.Range("M125").Value = 1.571
$ws.Range("N125").Value = 4
$ws.Range("O125").Value = 3.75
$ws.Range("P125").Value = 1.7
$ws.Range("Q125").Value = 0.75
$ws.Range("R125").Value = 1.875
$ws.Range("S125").Value = 1.925
$ws.Range("T125").Value = 2.5
$ws.Range("U125").Value = 1.825
$ws.Range("V125").Value = 1.975
$ws.Range("W125").Value = -1
$ws.Range("X125").Value = -1
$ws.Range("Y125").Value = 0.7
$ws.Range("Z125").Value = -1
$ws.Range("AA125").Value = 0.925
$ws.Range("AB125").Value = 0.825
$ws.Range("AC125").Value = -1

# ---------------------------------------------------------------------------
# Row 126: brand-new match row (upcoming - no result / closing odds yet)
# ---------------------------------------------------------------------------
$ws.Range("A124").Copy()
$ws.Range("A126").PasteSpecial(-4122)
$ws.Range("E124").Copy()
$ws.Range("E126").PasteSpecial(-4122)

$ws.Range("A126").Value = 124
$ws.Range("B126").Value = 7011620
$ws.Range("C126").Value = "Azerbaijan Premier League"
$ws.Range("D126").Value = "Azerbaijan Premier League"
$ws.Range("E126").Value = 45359.5
$ws.Range("F126").Value = "Sabail FC"
$ws.Range("G126").Value = "FK Sumqayit"
$ws.Range("K126").Value = 2.2
$ws.Range("L126").Value = 3.25
$ws.Range("M126").Value = 2.8
$ws.Range("N126").Value = 2.3
$ws.Range("O126").Value = 3.25
$ws.Range("P126").Value = 2.625
$ws.Range("Q126").Value = 0
$ws.Range("R126").Value = 1.775
$ws.Range("S126").Value = 2.025
$ws.Range("T126").Value = 2.25
$ws.Range("U126").Value = 1.975
$ws.Range("V126").Value = 1.825
$ws.Range("W126").Value = 0
$ws.Range("X126").Value = 0
$ws.Range("Y126").Value = 0
$ws.Range("Z126").Value = 0
$ws.Range("AA126").Value = 0
